$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 616.58826
$ws.Range("I53").Value = 415
$ws.Range("J53").Value = 726.5454999999999
$ws.Range("K53").Value = 415
$ws.Range("L53").Value = 726.5454999999999
$ws.Range("M53").Value = 222
$ws.Range("N53").Value = -2000.5455

# Row 111
$ws.Range("H111").Value = 758.1429000000001
$ws.Range("I111").Value = 639.25
$ws.Range("J111").Value = 916.6667
$ws.Range("K111").Value = 1917.75
$ws.Range("L111").Value = 2750.0001
$ws.Range("M111").Value = 1149.25
$ws.Range("N111").Value = -8884.000100000001

# Row 129
$ws.Range("H129").Value = 996.1277
$ws.Range("J129").Value = 1031.8667
$ws.Range("L129").Value = 3095.6001
$ws.Range("N129").Value = -13095.6001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7660.293
$ws.Range("I32").Value = 6244.4053
$ws.Range("K32").Value = 6244.4053
$ws.Range("M32").Value = -5957.4053

# Row 42
$ws.Range("H42").Value = 27999
$ws.Range("J42").Value = 27999
$ws.Range("L42").Value = 27999
$ws.Range("N42").Value = -28971

# Row 45
$ws.Range("H45").Value = 1160
$ws.Range("J45").Value = 1186.6666
$ws.Range("L45").Value = 1186.6666
$ws.Range("N45").Value = -1940.6666

# Row 110
$ws.Range("H110").Value = 1429.125
$ws.Range("I110").Value = 905.5
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 905.5
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 1139.5
$ws.Range("N110").Value = -7090

# Row 124
$ws.Range("H124").Value = 29389.25
$ws.Range("J124").Value = 29389.25
$ws.Range("L124").Value = 29389.25
$ws.Range("N124").Value = -39209.25

# Row 137
$ws.Range("H137").Value = 39591.668
$ws.Range("J137").Value = 39591.668
$ws.Range("L137").Value = 39591.668
$ws.Range("N137").Value = -49791.668

$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 81995.336
$ws.Range("J59").Value = 81995.336
$ws.Range("L59").Value = 81995.336
$ws.Range("N59").Value = -83689.336

# Row 86
$ws.Range("H86").Value = 1769.9
$ws.Range("I86").Value = 1649.875
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 1649.875
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -526.875
$ws.Range("N86").Value = -4496

# Row 89
$ws.Range("H89").Value = 1769.9
$ws.Range("I89").Value = 1649.875
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 8249.375
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -2633.375
$ws.Range("N89").Value = -22482

# Row 134
$ws.Range("H134").Value = 2415.9167
$ws.Range("I134").Value = 1370.8182
$ws.Range("J134").Value = 4715.1333
$ws.Range("K134").Value = 4112.4546
$ws.Range("L134").Value = 14145.3999
$ws.Range("M134").Value = -1577.4546
$ws.Range("N134").Value = -19215.3999

# Row 137
$ws.Range("H137").Value = 36195
$ws.Range("J137").Value = 39926.668
$ws.Range("L137").Value = 39926.668
$ws.Range("N137").Value = -50126.668

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1856.84
$ws.Range("I122").Value = 1347.75
$ws.Range("J122").Value = 2096.4119
$ws.Range("K122").Value = 4043.25
$ws.Range("L122").Value = 6289.2357
$ws.Range("M122").Value = -1593.25
$ws.Range("N122").Value = -11189.2357

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1877
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3384
$ws.Range("N84").ClearContents()

# Row 107
$ws.Range("H107").Value = 72103.42999999999
$ws.Range("I107").Value = 607.1667
$ws.Range("J107").Value = 125725.625
$ws.Range("K107").Value = 1821.5001
$ws.Range("L107").Value = 377176.875
$ws.Range("M107").Value = 98.49990000000003
$ws.Range("N107").Value = -381016.875

# Row 131
$ws.Range("H131").Value = 9434880
$ws.Range("I131").Value = 125000290
$ws.Range("J131").Value = 968.7755
$ws.Range("K131").Value = 375000870
$ws.Range("L131").Value = 2906.3265
$ws.Range("M131").Value = -374995830
$ws.Range("N131").Value = -12986.3265

# Row 132
$ws.Range("H132").Value = 1967.1666
$ws.Range("I132").Value = 799.25
$ws.Range("J132").Value = 2901.5
$ws.Range("K132").Value = 7193.25
$ws.Range("L132").Value = 26113.5
$ws.Range("M132").Value = -4663.25
$ws.Range("N132").Value = -31173.5

# Row 137
$ws.Range("H137").Value = 3634.3
$ws.Range("I137").Value = 626
$ws.Range("J137").Value = 4637.067
$ws.Range("K137").Value = 1878
$ws.Range("L137").Value = 13911.201
$ws.Range("M137").Value = 3222
$ws.Range("N137").Value = -24111.201

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5718.4727
$ws.Range("I70").Value = 5191.8423
$ws.Range("J70").Value = 6895.647
$ws.Range("K70").Value = 5191.8423
$ws.Range("L70").Value = 6895.647
$ws.Range("M70").Value = -4921.8423
$ws.Range("N70").Value = -7435.647

# Row 73
$ws.Range("H73").Value = 5718.4727
$ws.Range("I73").Value = 5191.8423
$ws.Range("J73").Value = 6895.647
$ws.Range("K73").Value = 5191.8423
$ws.Range("L73").Value = 6895.647
$ws.Range("M73").Value = -4255.8423
$ws.Range("N73").Value = -8767.647000000001

# Row 123
$ws.Range("H123").Value = 10980.315
$ws.Range("J123").Value = 10980.315
$ws.Range("L123").Value = 10980.315
$ws.Range("N123").Value = -15880.315

# Row 137
$ws.Range("H137").Value = 42398
$ws.Range("J137").Value = 42398
$ws.Range("L137").Value = 42398
$ws.Range("N137").Value = -52598

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3817.75
$ws.Range("I7").Value = 2069.5386
$ws.Range("J7").Value = 5332.8667
$ws.Range("K7").Value = 2069.5386
$ws.Range("L7").Value = 5332.8667
$ws.Range("M7").Value = -1957.5386
$ws.Range("N7").Value = -5556.8667

# Row 126
$ws.Range("H126").Value = 3817.75
$ws.Range("I126").Value = 2069.5386
$ws.Range("J126").Value = 5332.8667
$ws.Range("K126").Value = 6208.6158
$ws.Range("L126").Value = 15998.6001
$ws.Range("M126").Value = -3738.6158
$ws.Range("N126").Value = -20938.6001
